$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("C2").Value = 1.032778917184339
$ws.Range("D2").Value = 1.04094639782437
$ws.Range("E2").Value = 1.041521604215282
$ws.Range("F2").Value = 1.050974861329929
$ws.Range("J2").Value = 1.037907263792015
$ws.Range("K2").Value = 1.043727305393431
$ws.Range("L2").Value = 1.044300882976482
$ws.Range("M2").Value = 1.053727648330374
$ws.Range("N2").Value = 1.016511895717824
$ws.Range("C3").Value = 1.033887469934978
$ws.Range("D3").Value = 1.041948041778461
$ws.Range("E3").Value = 1.042516183855933
$ws.Range("F3").Value = 1.052073023416918
$ws.Range("J3").Value = 1.038657446785932
$ws.Range("K3").Value = 1.044538974358905
$ws.Range("L3").Value = 1.045105624527215
$ws.Range("M3").Value = 1.054637628989984
$ws.Range("N3").Value = 1.016769671234599
$ws.Range("C4").Value = 1.034605539528362
$ws.Range("D4").Value = 1.042597193449364
$ws.Range("E4").Value = 1.043160767956477
$ws.Range("F4").Value = 1.052784864318794
$ws.Range("J4").Value = 1.039143086872412
$ws.Range("K4").Value = 1.045064588212531
$ws.Range("L4").Value = 1.045626757672421
$ws.Range("M4").Value = 1.055227109970175
$ws.Range("N4").Value = 1.016936320762673
$ws.Range("C5").Value = 1.034907598564819
$ws.Range("D5").Value = 1.042870340833675
$ws.Range("E5").Value = 1.043431995946438
$ws.Range("F5").Value = 1.053084422705107
$ws.Range("J5").Value = 1.039347302746772
$ws.Range("K5").Value = 1.045285654064283
$ws.Range("L5").Value = 1.045845940244149
$ws.Range("M5").Value = 1.055475086085017
$ws.Range("N5").Value = 1.017006344425932
$ws.Range("C6").Value = 1.034958326318581
$ws.Range("D6").Value = 1.042916217778877
$ws.Range("E6").Value = 1.043477550663671
$ws.Range("F6").Value = 1.053134737485721
$ws.Range("J6").Value = 1.03938159457864
$ws.Range("K6").Value = 1.045322777718984
$ws.Range("L6").Value = 1.045882747711316
$ws.Range("M6").Value = 1.055516731639881
$ws.Range("N6").Value = 1.017018099600504
$ws.Range("C7").Value = 1.034609574938314
$ws.Range("D7").Value = 1.042600842299759
$ws.Range("E7").Value = 1.043164391156677
$ws.Range("F7").Value = 1.052788865852406
$ws.Range("J7").Value = 1.03914581540783
$ws.Range("K7").Value = 1.045067541721402
$ws.Range("L7").Value = 1.045629686015226
$ws.Range("M7").Value = 1.055230422817505
$ws.Range("N7").Value = 1.016937256563269
$ws.Range("C8").Value = 1.033153400147239
$ws.Range("D8").Value = 1.041284696417669
$ws.Range("E8").Value = 1.041857514646189
$ws.Range("F8").Value = 1.051345729942387
$ws.Range("J8").Value = 1.038160745625785
$ws.Range("K8").Value = 1.044001527689888
$ws.Range("L8").Value = 1.044572763738085
$ws.Range("M8").Value = 1.054035043151074
$ws.Range("N8").Value = 1.016599042748458
$ws.Range("C9").Value = 1.030593260206127
$ws.Range("D9").Value = 1.038973317793849
$ws.Range("E9").Value = 1.039562498016328
$ws.Range("F9").Value = 1.048812380038184
$ws.Range("J9").Value = 1.036426634509789
$ws.Range("K9").Value = 1.042126226636269
$ws.Range("L9").Value = 1.042713498095249
$ws.Range("M9").Value = 1.051933720943188
$ws.Range("N9").Value = 1.016001940929356
$ws.Range("C10").Value = 1.028890403018466
$ws.Range("D10").Value = 1.037437692101012
$ws.Range("E10").Value = 1.038037801025488
$ws.Range("F10").Value = 1.047129980899863
$ws.Range("J10").Value = 1.035271721961216
$ws.Range("K10").Value = 1.040878161823052
$ws.Range("L10").Value = 1.041476134541979
$ws.Range("M10").Value = 1.050536279062444
$ws.Range("N10").Value = 1.015603130729733
$ws.Range("C11").Value = 1.028153968821848
$ws.Range("D11").Value = 1.03677400903684
$ws.Range("E11").Value = 1.037378855507853
$ws.Range("F11").Value = 1.046403028840856
$ws.Range("J11").Value = 1.034771909872483
$ws.Range("K11").Value = 1.04033824558155
$ws.Range("L11").Value = 1.040940854931738
$ws.Range("M11").Value = 1.049931990277615
$ws.Range("N11").Value = 1.015430268366702
$ws.Range("C12").Value = 1.027880561293764
$ws.Range("D12").Value = 1.036527676121872
$ws.Range("E12").Value = 1.037134283165018
$ws.Range("F12").Value = 1.046133237728194
$ws.Range("J12").Value = 1.034586298490851
$ws.Range("K12").Value = 1.040137772552843
$ws.Range("L12").Value = 1.040742104614951
$ws.Range("M12").Value = 1.04970765280972
$ws.Range("N12").Value = 1.015366033521727
$ws.Range("C13").Value = 1.027939201947364
$ws.Range("D13").Value = 1.036580506832641
$ws.Range("E13").Value = 1.037186736188707
$ws.Range("F13").Value = 1.046191098369845
$ws.Range("J13").Value = 1.034626110901285
$ws.Range("K13").Value = 1.040180771260392
$ws.Range("L13").Value = 1.040784733773739
$ws.Range("M13").Value = 1.049755768421109
$ws.Range("N13").Value = 1.015379813290175
$ws.Range("C14").Value = 1.028131366079542
$ws.Range("D14").Value = 1.036753643233756
$ws.Range("E14").Value = 1.037358635216957
$ws.Range("F14").Value = 1.046380723104707
$ws.Range("J14").Value = 1.034756566339011
$ws.Range("K14").Value = 1.040321672869031
$ws.Range("L14").Value = 1.040924424609596
$ws.Range("M14").Value = 1.04991344397672
$ws.Range("N14").Value = 1.015424959224897
$ws.Range("C15").Value = 1.02824978287455
$ws.Range("D15").Value = 1.036860343274568
$ws.Range("E15").Value = 1.037464572983329
$ws.Range("F15").Value = 1.046497587811239
$ws.Range("J15").Value = 1.034836949687172
$ws.Range("K15").Value = 1.040408497063664
$ws.Range("L15").Value = 1.041010502868871
$ws.Range("M15").Value = 1.050010609301043
$ws.Range("N15").Value = 1.015452771677571
$ws.Range("C16").Value = 1.028939296955285
$ws.Range("D16").Value = 1.037481764982884
$ws.Range("E16").Value = 1.038081559607479
$ws.Range("F16").Value = 1.047178258732104
$ws.Range("J16").Value = 1.035304898598095
$ws.Range("K16").Value = 1.040914004920491
$ws.Range("L16").Value = 1.041511669984828
$ws.Range("M16").Value = 1.050576400858557
$ws.Range("N16").Value = 1.015614599368302
$ws.Range("C17").Value = 1.029372055375552
$ws.Range("D17").Value = 1.037871902059073
$ws.Range("E17").Value = 1.038468916127405
$ws.Range("F17").Value = 1.047605637489823
$ws.Range("J17").Value = 1.03559850366848
$ws.Range("K17").Value = 1.041231231715783
$ws.Range("L17").Value = 1.041826174744545
$ws.Range("M17").Value = 1.050931524683973
$ws.Range("N17").Value = 1.01571606292318
$ws.Range("C18").Value = 1.029624564275757
$ws.Range("D18").Value = 1.038099583291894
$ws.Range("E18").Value = 1.038694976022443
$ws.Range("F18").Value = 1.047855068656765
$ws.Range("J18").Value = 1.035769784838769
$ws.Range("K18").Value = 1.041416313360777
$ws.Range("L18").Value = 1.042009668951992
$ws.Range("M18").Value = 1.051138740811138
$ws.Range("N18").Value = 1.015775228007061
$ws.Range("C19").Value = 1.029710678264679
$ws.Range("D19").Value = 1.038177237202941
$ws.Range("E19").Value = 1.038772077185316
$ws.Range("F19").Value = 1.047940143398595
$ws.Range("J19").Value = 1.035828191749015
$ws.Range("K19").Value = 1.041479429660769
$ws.Range("L19").Value = 1.042072244023563
$ws.Range("M19").Value = 1.05120940944944
$ws.Range("N19").Value = 1.015795398900244
$ws.Range("C20").Value = 1.029325615369017
$ws.Range("D20").Value = 1.037830031538081
$ws.Range("E20").Value = 1.03842734388042
$ws.Range("F20").Value = 1.04755976844583
$ws.Range("J20").Value = 1.03556699990182
$ws.Range("K20").Value = 1.041197191252291
$ws.Range("L20").Value = 1.041792426300611
$ws.Range("M20").Value = 1.050893415146861
$ws.Range("N20").Value = 1.01570517858955
$ws.Range("C21").Value = 1.028074774755831
$ws.Range("D21").Value = 1.036702653672456
$ws.Range("E21").Value = 1.037308010007771
$ws.Range("F21").Value = 1.046324876961287
$ws.Range("J21").Value = 1.034718149329076
$ws.Range("K21").Value = 1.040280178760443
$ws.Range("L21").Value = 1.040883287030518
$ws.Range("M21").Value = 1.049867009078174
$ws.Range("N21").Value = 1.015411665590211
$ws.Range("C22").Value = 1.027289113602919
$ws.Range("D22").Value = 1.035994917354999
$ws.Range("E22").Value = 1.036605336123457
$ws.Range("F22").Value = 1.045549789715973
$ws.Range("J22").Value = 1.034184681034618
$ws.Range("K22").Value = 1.039704055853221
$ws.Range("L22").Value = 1.040312117006445
$ws.Range("M22").Value = 1.049222374488054
$ws.Range("N22").Value = 1.015226971610363
$ws.Range("C23").Value = 1.027705532513928
$ws.Range("D23").Value = 1.036369998269675
$ws.Range("E23").Value = 1.036977732875266
$ws.Range("F23").Value = 1.045960551229606
$ws.Range("J23").Value = 1.034467460094397
$ws.Range("K23").Value = 1.04000942781228
$ws.Range("L23").Value = 1.040614863087236
$ws.Range("M23").Value = 1.049564040342087
$ws.Range("N23").Value = 1.015324895597774
$ws.Range("C24").Value = 1.029346599336147
$ws.Range("D24").Value = 1.037848950650665
$ws.Range("E24").Value = 1.03844612821424
$ws.Range("F24").Value = 1.04758049423388
$ws.Range("J24").Value = 1.035581235016567
$ws.Range("K24").Value = 1.041212572522179
$ws.Range("L24").Value = 1.041807675618571
$ws.Range("M24").Value = 1.050910634963929
$ws.Range("N24").Value = 1.015710096803594
$ws.Range("C25").Value = 1.031254428840191
$ws.Range("D25").Value = 1.039569933292317
$ws.Range("E25").Value = 1.040154880025949
$ws.Range("F25").Value = 1.04946616734946
$ws.Range("J25").Value = 1.036874740070764
$ws.Range("K25").Value = 1.042610661690447
$ws.Range("L25").Value = 1.043193785761264
$ws.Range("M25").Value = 1.052476358168659
$ws.Range("N25").Value = 1.016156437977594
